$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet used to have 5 columns:
#   A: Display Type | B: Product | C: Description | D: Quantity | E: Unit Price
# The edit drops the old "Display Type" column and shifts everything one
# column to the left, renaming the old "Product" header to "Display".
# (Row 3's "Note" cell keeps its own original formatting - that's the one
# cell that does NOT pick up its right-hand neighbour's style.)
# ---------------------------------------------------------------------------

# --- 1) Copy formatting from the old layout into the new column positions.
# Do this BEFORE touching any values, since PasteSpecial(formats) only reads
# the current formatting of the source cell (not its value), so order here
# does not corrupt any of the source data we still need to read from.

function CopyFmt($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

CopyFmt "B1" "A1"
CopyFmt "C1" "B1"
CopyFmt "D1" "C1"
CopyFmt "E1" "D1"

CopyFmt "B2" "A2"
CopyFmt "C2" "B2"
CopyFmt "D2" "C2"
CopyFmt "E2" "D2"

# A3 keeps its own existing formatting - no-op.
CopyFmt "C3" "B3"
CopyFmt "D3" "C3"
CopyFmt "E3" "D3"

CopyFmt "B4" "A4"
CopyFmt "C4" "B4"
CopyFmt "D4" "C4"
CopyFmt "E4" "D4"

CopyFmt "B5" "A5"
CopyFmt "C5" "B5"
CopyFmt "D5" "C5"
CopyFmt "E5" "D5"

$excel.CutCopyMode = 0

# --- 2) Now write the actual values/content for the new A:D layout.
$ws.Range("A1").Value = "Display"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Quantity"
$ws.Range("D1").Value = "Unit Price"

$ws.Range("A2").Value = $null
$ws.Range("B2").Value = "Description row 2"
$ws.Range("C2").Value = $null
$ws.Range("D2").Value = $null

$ws.Range("A3").Value = "Note"
$ws.Range("B3").Value = "Description row 3"
$ws.Range("C3").Value = $null
$ws.Range("D3").Value = $null

$ws.Range("A4").Value = "Server Config"
$ws.Range("B4").Value = "Setup Server 1"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1600

$ws.Range("A5").Value = "Server Config"
$ws.Range("B5").Value = "Setup Server 2"
$ws.Range("C5").Value = $null
$ws.Range("D5").Value = $null

# --- 3) Drop the now-unused rightmost column (old column E).
$ws.Columns("E").Clear()

# --- 4) Row heights: rows 2-5 go back to the sheet's default (no explicit
# row height override any more); row 1 keeps its existing explicit height.
$ws.Rows("2:5").AutoFit()

# --- 5) Column widths for the new A:D columns (matching the widths that
# used to belong to old columns B:E).
$ws.Columns("A").ColumnWidth = 11
$ws.Columns("B").ColumnWidth = 35.17
$ws.Columns("C").ColumnWidth = 7
$ws.Columns("D").ColumnWidth = 4.17

# --- 6) Conditional formatting: shift every rule one column to the left,
# and point every formula at column A instead of column B.
$cf1 = $ws.Range("B2:B5").FormatConditions.Item(1)
$cf1.Formula1 = "=LEN(TRIM(A2))>0"
$cf1.ModifyAppliesToRange($ws.Range("A2:A5")) | Out-Null

$cf2 = $ws.Range("C2:C5").FormatConditions.Item(1)
$cf2.Formula1 = "=NOT(ISBLANK(A2))"
$cf2.ModifyAppliesToRange($ws.Range("B2:B5")) | Out-Null

$cf3 = $ws.Range("D2:D5").FormatConditions.Item(1)
$cf3.Formula1 = "=NOT(ISBLANK(A2))"
$cf3.ModifyAppliesToRange($ws.Range("C2:C5")) | Out-Null

$cf4 = $ws.Range("E2:E5").FormatConditions.Item(1)
$cf4.Formula1 = "=NOT(ISBLANK(A2))"
$cf4.ModifyAppliesToRange($ws.Range("D2:D5")) | Out-Null

# --- 7) Selection moves to A5.
$ws.Range("A5").Select()
